# Updates the cryptocurrency price/volume table on Sheet1 to the latest
# scraped snapshot (GitHub Actions run). Column D = Price, Column E = Volume(1h).
# A handful of rows also changed identity (coin name/link) because the
# source ranking reordered two adjacent entries (Dai/Litecoin at rows 24-25,
# BabyDogeCoin/Aave at rows 46-47).
# Price values are stored as text (leading apostrophe keeps Excel from
# auto-converting plain-number-looking strings like "591.32" into numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.700.53'
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").Value = '2.513.22'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'" + '591.32'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = "'" + '172.61'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("E8").Value = '  -0.65%  '
$ws.Range("D9").Value = '2.511.51'
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("D12").Value = "'" + '5.12'
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("D13").Value = "'" + '0.341'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").Value = "'" + '26.39'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '2.944.71'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '67.516.83'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = '2.468.14'
$ws.Range("E18").Value = '  -3.56%  '
$ws.Range("D19").Value = "'" + '11.82'
$ws.Range("E19").Value = '  +4.64%  '
$ws.Range("D20").Value = "'" + '7.87'
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("D21").Value = "'" + '366.99'
$ws.Range("E21").Value = '  +3.69%  '
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("D23").Value = "'" + '4.58'
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = "'" + '1.01'
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = "'" + '71.42'
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").Value = "'" + '1.92'
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("D27").Value = "'" + '9.97'
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").Value = "'" + '0.997'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").Value = '2.601.91'
$ws.Range("E29").Value = '  -2.89%  '
$ws.Range("E30").Value = '  -2.32%  '
$ws.Range("D31").Value = "'" + '8.36'
$ws.Range("E31").Value = '  +2.68%  '
$ws.Range("D32").Value = "'" + '534.76'
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("E33").Value = '  -1.65%  '
$ws.Range("E34").Value = '  +1.30%  '
$ws.Range("E35").Value = '  -2.47%  '
$ws.Range("D36").Value = "'" + '1.00'
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").Value = "'" + '158.52'
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("D38").Value = "'" + '1.43'
$ws.Range("E38").Value = '  -2.37%  '
$ws.Range("D39").Value = "'" + '18.94'
$ws.Range("E39").Value = '  +1.82%  '
$ws.Range("E40").Value = '  +1.03%  '
$ws.Range("D41").Value = "'" + '0.350'
$ws.Range("E41").Value = '  -1.55%  '
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").Value = "'" + '5.13'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").Value = "'" + '2.47'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₆0281'
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = "'" + '146.13'
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").Value = "'" + '3.70'
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("D49").Value = "'" + '0.550'
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("E50").Value = '  +1.38%  '
$ws.Range("D51").Value = "'" + '0.0750'
$ws.Range("E51").Value = '  -1.35%  '
